$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the soil-class annotation text that was in column K for rows 13-25
$rows = 13..25
foreach ($r in $rows) {
    $ws.Cells.Item($r, 11).Value = $null
}

# Reset custom row heights back to the sheet default for the affected rows
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(16).AutoFit()
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(20).AutoFit()
$ws.Rows.Item(21).AutoFit()

# Update the selected range shown when the sheet is active
$ws.Range("L1:N1").Select()
